# Update "want-to-go" / "interested" counts (column F) on several sheets,
# matching the regenerated data snapshot described in the commit message
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F8").Value  = 8083
$wsExpo.Range("F14").Value = 4939
$wsExpo.Range("F17").Value = 5356
$wsExpo.Range("F20").Value = 328
$wsExpo.Range("F21").Value = 449
$wsExpo.Range("F27").Value = 9078
$wsExpo.Range("F29").Value = 1632
$wsExpo.Range("F36").Value = 1006
$wsExpo.Range("F37").Value = 1177
$wsExpo.Range("F39").Value = 4730
$wsExpo.Range("F41").Value = 375
$wsExpo.Range("F42").Value = 1158

# Sheet "演出" (Performances)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F9").Value = 179

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 8083
$wsAll.Range("F17").Value = 4939
$wsAll.Range("F19").Value = 5356
$wsAll.Range("F22").Value = 328
$wsAll.Range("F23").Value = 449
$wsAll.Range("F29").Value = 179
$wsAll.Range("F30").Value = 9078
$wsAll.Range("F32").Value = 1632
$wsAll.Range("F38").Value = 1006
$wsAll.Range("F39").Value = 1177
$wsAll.Range("F41").Value = 4730
$wsAll.Range("F42").Value = 375
$wsAll.Range("F43").Value = 1158
